$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.100.37"
$ws.Range("E2").Value = "  +1.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.190.69"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'616.03"
$ws.Range("E5").Value = "  +1.57%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'147.79"

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.189.60"
$ws.Range("E8").Value = "  +1.39%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +0.11%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.42%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'5.53"
$ws.Range("E11").Value = "  -1.14%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.21%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  +1.63%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'36.19"
$ws.Range("E14").Value = "  -2.03%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.713.88"
$ws.Range("E15").Value = "  +1.48%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +3.14%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "65.121.08"
$ws.Range("E17").Value = "  +1.17%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.182.52"
$ws.Range("E18").Value = "  +1.40%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.97"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'485.45"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'14.82"
$ws.Range("E21").Value = "  +1.11%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'8.01"
$ws.Range("E23").Value = "  +2.82%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'14.04"
$ws.Range("E24").Value = "  +1.31%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'84.82"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.02%  "

# Row 27 - was PancakeSwap, becomes RenderToken
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'8.76"
$ws.Range("E27").Value = "  +1.95%  "

# Row 28 - was RenderToken, becomes PancakeSwap
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.84"
$ws.Range("E28").Value = "  -3.39%  "

# Row 29 - NEARProtocol
$ws.Range("D29").Value = "'7.14"
$ws.Range("E29").Value = "  +1.23%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -3.66%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "'2.14"
$ws.Range("E31").Value = "  -4.94%  "

# Row 32 - Stacks
$ws.Range("D32").Value = "'2.74"
$ws.Range("E32").Value = "  +0.50%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.15%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'26.90"
$ws.Range("E34").Value = "  +0.48%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +2.64%  "

# Row 36 - PEPE
$ws.Range("D36").Value = "0.0₃0801"
$ws.Range("E36").Value = "  +5.72%  "

# Row 37 - Filecoin
$ws.Range("D37").Value = "'6.08"
$ws.Range("E37").Value = "  -0.59%  "

# Row 38 - dogwifhat
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  -1.25%  "

# Row 39 - OKB
$ws.Range("D39").Value = "'53.35"
$ws.Range("E39").Value = "  -2.09%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "'471.19"
$ws.Range("E40").Value = "  +3.82%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.0404"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -2.35%  "

# Row 43 - Cosmos
$ws.Range("D43").Value = "'8.45"
$ws.Range("E43").Value = "  -0.69%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.876.42"
$ws.Range("E44").Value = "  -0.48%  "

# Row 45 - Fetch.AI
$ws.Range("D45").Value = "'2.36"
$ws.Range("E45").Value = "  +1.56%  "

# Row 46 - TheGraph
$ws.Range("D46").Value = "'0.273"
$ws.Range("E46").Value = "  -0.38%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  +6.40%  "

# Row 48 - Arweave
$ws.Range("D48").Value = "'37.30"
$ws.Range("E48").Value = "  +10.52%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'27.05"
$ws.Range("E49").Value = "  +0.94%  "

# Row 50 - USDe
$ws.Range("E50").Value = "  +0.12%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.43%  "
